$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Set width of new column G (OOXML stored width 17 == ColumnWidth property 16.17 for this font)
$ws.Columns.Item(7).ColumnWidth = 16.17

# Copy formatting from column F into the new column G (header style, body style, totals style)
$ws.Range("F1:F56").Copy()
$ws.Range("G1:G56").PasteSpecial(-4122)  # xlPasteFormats

# Header text for the new column
$ws.Range("G1").Value = "PRESUPUESTO"

# Data + totals rows all contain 0 in the new column
$ws.Range("G2:G56").Value = 0
